$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price/volume figures (cryptos list refresh).
# Numeric-looking Price values are prefixed with an apostrophe so Excel
# stores them as text (matching the original inlineStr cell content)
# instead of silently converting them to numbers.

$ws.Range('D2').Value = '51.940.58'
$ws.Range('E2').Value = '  +0.43%  '
$ws.Range('D3').Value = '2.936.78'
$ws.Range('E3').Value = '  +3.86%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '''352.43'
$ws.Range('E5').Value = '  +0.63%  '
$ws.Range('D6').Value = '''112.35'
$ws.Range('E6').Value = '  -0.29%  '
$ws.Range('D7').Value = '''0.560'
$ws.Range('E7').Value = '  +0.49%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = '''0.628'
$ws.Range('E9').Value = '  +1.68%  '
$ws.Range('D10').Value = '''39.42'
$ws.Range('E10').Value = '  -1.66%  '
$ws.Range('D11').Value = '''0.0889'
$ws.Range('E11').Value = '  +5.02%  '
$ws.Range('D13').Value = '''20.01'
$ws.Range('E13').Value = '  +0.50%  '
$ws.Range('D14').Value = '''7.84'
$ws.Range('E14').Value = '  +1.10%  '
$ws.Range('D15').Value = '3.397.82'
$ws.Range('E15').Value = '  +3.92%  '
$ws.Range('D16').Value = '2.933.33'
$ws.Range('E16').Value = '  +4.11%  '
$ws.Range('D17').Value = '''0.991'
$ws.Range('E17').Value = '  +0.92%  '
$ws.Range('D18').Value = '51.992.95'
$ws.Range('E18').Value = '  +0.49%  '
$ws.Range('D19').Value = '''7.65'
$ws.Range('E19').Value = '  +0.84%  '
$ws.Range('D20').Value = '''3.32'
$ws.Range('E20').Value = '  -3.78%  '
$ws.Range('D21').Value = '''14.27'
$ws.Range('E21').Value = '  +6.82%  '
$ws.Range('D22').Value = '0.0₃0987'
$ws.Range('E22').Value = '  +1.63%  '
$ws.Range('D23').Value = '''71.26'
$ws.Range('E23').Value = '  +1.27%  '
$ws.Range('D24').Value = '''269.35'
$ws.Range('E24').Value = '  +0.18%  '
$ws.Range('E25').Value = '  +1.40%  '
$ws.Range('E26').Value = '  +9.91%  '
$ws.Range('D27').Value = '''26.98'
$ws.Range('E27').Value = '  +2.86%  '
$ws.Range('E28').Value = '  +0.15%  '
$ws.Range('D29').Value = '''7.43'
$ws.Range('E29').Value = '  +16.93%  '
$ws.Range('E30').Value = '  +20.84%  '
$ws.Range('D31').Value = '''10.61'
$ws.Range('E31').Value = '  +0.54%  '
$ws.Range('D32').Value = '''37.54'
$ws.Range('E32').Value = '  -2.40%  '
$ws.Range('E33').Value = '  +0.21%  '
$ws.Range('E34').Value = '  +10.46%  '
$ws.Range('D35').Value = '''52.92'
$ws.Range('E35').Value = '  +0.23%  '
$ws.Range('D36').Value = '''0.0454'
$ws.Range('E36').Value = '  +1.38%  '
$ws.Range('D37').Value = '''0.998'
$ws.Range('E37').Value = '  -0.15%  '
$ws.Range('D38').Value = '''3.30'
$ws.Range('E38').Value = '  +3.03%  '
$ws.Range('E39').Value = '  -0.13%  '
$ws.Range('E40').Value = '  +1.97%  '
$ws.Range('D41').Value = '''2.71'
$ws.Range('E41').Value = '  +7.71%  '
$ws.Range('E42').Value = '  +1.69%  '
$ws.Range('D43').Value = '''23.17'
$ws.Range('E43').Value = '  +5.15%  '
$ws.Range('E44').Value = '  -0.72%  '
$ws.Range('D46').Value = '''3.54'
$ws.Range('E46').Value = '  +1.46%  '
$ws.Range('D47').Value = '2.171.38'
$ws.Range('E47').Value = '  +0.14%  '
$ws.Range('D48').Value = '''112.23'
$ws.Range('E48').Value = '  -8.25%  '
$ws.Range('D49').Value = '''0.249'
$ws.Range('E49').Value = '  +1.40%  '
$ws.Range('E50').Value = '  +10.86%  '
$ws.Range('E51').Value = '  -1.56%  '
